# Apply latest crypto price/volume snapshot to Sheet1 (D = Price, E = Volume(1h)).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.642.99"
$ws.Range("E2").Value = "  +0.49%  "
$ws.Range("D3").Value = "1.813.48"
$ws.Range("E3").Value = "  +0.57%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.24%  "
$ws.Range("D5").Value = "'226.20"
$ws.Range("E5").Value = "  -1.04%  "
$ws.Range("E6").Value = "  +3.58%  "
$ws.Range("E7").Value = "  -0.25%  "
$ws.Range("D8").Value = "'37.06"
$ws.Range("E8").Value = "  +6.36%  "
$ws.Range("D9").Value = "'0.294"
$ws.Range("E9").Value = "  -2.21%  "
$ws.Range("D10").Value = "'0.0685"
$ws.Range("E10").Value = "  -1.35%  "
$ws.Range("E11").Value = "  +1.50%  "
$ws.Range("D12").Value = "2.073.93"
$ws.Range("E12").Value = "  +0.51%  "
$ws.Range("D13").Value = "'11.38"
$ws.Range("E13").Value = "  +1.48%  "
$ws.Range("D14").Value = "1.834.74"
$ws.Range("E14").Value = "  +1.55%  "
$ws.Range("D15").Value = "'0.635"
$ws.Range("E15").Value = "  -1.32%  "
$ws.Range("D16").Value = "34.588.37"
$ws.Range("E16").Value = "  +0.38%  "
$ws.Range("E17").Value = "  +2.24%  "
$ws.Range("D18").Value = "'68.90"
$ws.Range("E18").Value = "  -0.22%  "
$ws.Range("D19").Value = "'243.93"
$ws.Range("E19").Value = "  -0.57%  "
$ws.Range("D20").Value = "0.0₃0780"
$ws.Range("E20").Value = "  -2.24%  "
$ws.Range("D21").Value = "'11.27"
$ws.Range("E21").Value = "  -2.08%  "
$ws.Range("E23").Value = "  -1.01%  "
$ws.Range("D24").Value = "'2.22"
$ws.Range("E24").Value = "  +4.57%  "
$ws.Range("D25").Value = "'171.89"
$ws.Range("E25").Value = "  -1.11%  "
$ws.Range("D26").Value = "'7.91"
$ws.Range("E26").Value = "  +1.13%  "
$ws.Range("D27").Value = "'17.29"
$ws.Range("E27").Value = "  +2.78%  "
$ws.Range("E28").Value = "  +2.08%  "
$ws.Range("D29").Value = "'0.999"
$ws.Range("E29").Value = "  -0.25%  "
$ws.Range("E30").Value = "  +0.23%  "
$ws.Range("D31").Value = "'3.94"
$ws.Range("E31").Value = "  -1.93%  "
$ws.Range("D32").Value = "'1.24"
$ws.Range("E32").Value = "  -0.63%  "
$ws.Range("D33").Value = "'0.0518"
$ws.Range("E33").Value = "  -2.49%  "
$ws.Range("E34").Value = "  -0.58%  "
$ws.Range("D35").Value = "1.367.31"
$ws.Range("E35").Value = "  -1.96%  "
$ws.Range("D36").Value = "'0.656"
$ws.Range("E36").Value = "  -3.96%  "
$ws.Range("D37").Value = "'1.06"
$ws.Range("E37").Value = "  +0.47%  "
$ws.Range("D38").Value = "'2.38"
$ws.Range("E38").Value = "  -4.74%  "
$ws.Range("D39").Value = "'0.0188"
$ws.Range("E39").Value = "  -1.37%  "
$ws.Range("E40").Value = "  +0.82%  "
$ws.Range("D41").Value = "'81.56"
$ws.Range("E41").Value = "  -2.19%  "
$ws.Range("E42").Value = "  -1.53%  "
$ws.Range("D43").Value = "'0.942"
$ws.Range("E44").Value = "  +5.04%  "
$ws.Range("D45").Value = "'13.75"
$ws.Range("E45").Value = "  +0.83%  "
$ws.Range("E46").Value = "  -1.02%  "
$ws.Range("D47").Value = "1.974.26"
$ws.Range("E47").Value = "  +0.55%  "
$ws.Range("D48").Value = "'5.84"
$ws.Range("E48").Value = "  -2.23%  "
$ws.Range("D49").Value = "'0.999"
$ws.Range("E49").Value = "  -0.25%  "
$ws.Range("D50").Value = "'103.11"
$ws.Range("E50").Value = "  -1.79%  "
$ws.Range("E51").Value = "  -6.00%  "
